$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Cells whose new values look like plain numbers need to be forced to
# Text format first, otherwise Excel will silently convert them to
# numeric values (losing formatting / introducing float rounding).
$textCells = @('D5', 'D9', 'D10', 'D15', 'D17', 'D18', 'D20', 'D23', 'D24', 'D25', 'D26', 'D27', 'D30', 'D36', 'D37', 'D39', 'D42', 'D45', 'D48')
foreach ($cellAddr in $textCells) {
    $ws.Range($cellAddr).NumberFormat = "@"
}

$ws.Range('D2').Value = '26.874.04'
$ws.Range('E2').Value = '  +0.03%  '
$ws.Range('D3').Value = '1.543.38'
$ws.Range('E3').Value = '  -1.23%  '
$ws.Range('E4').Value = '  +0.34%  '
$ws.Range('D5').Value = '205.86'
$ws.Range('E5').Value = '  -0.09%  '
$ws.Range('E6').Value = '  -0.50%  '
$ws.Range('E7').Value = '  +0.29%  '
$ws.Range('E8').Value = '  -0.09%  '
$ws.Range('D9').Value = '21.25'
$ws.Range('D10').Value = '0.0581'
$ws.Range('E10').Value = '  -0.36%  '
$ws.Range('E11').Value = '  -0.75%  '
$ws.Range('D12').Value = '1.763.78'
$ws.Range('E12').Value = '  -1.18%  '
$ws.Range('D13').Value = '1.542.18'
$ws.Range('E13').Value = '  -1.31%  '
$ws.Range('E14').Value = '  -0.95%  '
$ws.Range('D15').Value = '0.511'
$ws.Range('E15').Value = '  -0.65%  '
$ws.Range('D16').Value = '26.863.78'
$ws.Range('E16').Value = '  -0.01%  '
$ws.Range('D17').Value = '61.34'
$ws.Range('E17').Value = '  +0.11%  '
$ws.Range('D18').Value = '213.58'
$ws.Range('E18').Value = '  -0.48%  '
$ws.Range('D19').Value = '0.0₃0680'
$ws.Range('E19').Value = '  +0.14%  '
$ws.Range('D20').Value = '7.15'
$ws.Range('E20').Value = '  -2.74%  '
$ws.Range('E21').Value = '  +0.25%  '
$ws.Range('E22').Value = '  -2.78%  '
$ws.Range('D23').Value = '9.14'
$ws.Range('E23').Value = '  -0.18%  '
$ws.Range('D24').Value = '1.93'
$ws.Range('E24').Value = '  -3.27%  '
$ws.Range('D25').Value = '152.27'
$ws.Range('E25').Value = '  -1.07%  '
$ws.Range('D26').Value = '6.58'
$ws.Range('E26').Value = '  -2.26%  '
$ws.Range('D27').Value = '14.78'
$ws.Range('E27').Value = '  -0.88%  '
$ws.Range('E28').Value = '  +0.28%  '
$ws.Range('E29').Value = '  -0.10%  '
$ws.Range('D30').Value = '0.0457'
$ws.Range('E30').Value = '  -1.25%  '
$ws.Range('E31').Value = '  -0.88%  '
$ws.Range('E32').Value = '  +1.55%  '
$ws.Range('D33').Value = '1.354.89'
$ws.Range('E33').Value = '  -3.41%  '
$ws.Range('E34').Value = '  +0.41%  '
$ws.Range('E35').Value = '  +0.13%  '
$ws.Range('D36').Value = '0.961'
$ws.Range('E36').Value = '  +4.83%  '
$ws.Range('D37').Value = '2.27'
$ws.Range('E37').Value = '  +0.26%  '
$ws.Range('E38').Value = '  -0.27%  '
$ws.Range('D39').Value = '0.518'
$ws.Range('E39').Value = '  -1.54%  '
$ws.Range('E40').Value = '  -1.09%  '
$ws.Range('E41').Value = '  +0.25%  '
$ws.Range('D42').Value = '5.58'
$ws.Range('E42').Value = '  +3.46%  '
$ws.Range('E43').Value = '  -0.79%  '
$ws.Range('E44').Value = '  +1.98%  '
$ws.Range('D45').Value = '63.18'
$ws.Range('E45').Value = '  -0.01%  '
$ws.Range('E46').Value = '  -1.62%  '
$ws.Range('D47').Value = '1.678.03'
$ws.Range('E47').Value = '  -1.27%  '
$ws.Range('D48').Value = '85.60'
$ws.Range('E48').Value = '  -0.68%  '
$ws.Range('E49').Value = '  +1.16%  '
$ws.Range('D50').Value = '0.0₇0978'
$ws.Range('E50').Value = '  +0.29%  '
$ws.Range('E51').Value = '  -0.04%  '
